# Update violent-crime-full-year workbook with 2024-05-28 data.
# Generated from the OOXML diff: for each affected sheet/cell, set the
# new 2024 (and a couple of 2023 correction) values, verifying the prior
# value first so a mismatch fails loudly instead of silently.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{Sheet="Citywide Totals"; Cells=@(@{Cell="K2"; Old=3007; New=3033}, @{Cell="J3"; Old=8077; New=8078}, @{Cell="K3"; Old=2961; New=2983}, @{Cell="J4"; Old=1818; New=1817}, @{Cell="K4"; Old=606; New=611}, @{Cell="K5"; Old=192; New=196}, @{Cell="K6"; Old=3549; New=3565}, @{Cell="K7"; Old=10315; New=10388})},
    @{Sheet="Logan Square"; Cells=@(@{Cell="K2"; Old=31; New=32}, @{Cell="K6"; Old=74; New=75}, @{Cell="K7"; Old=142; New=144})},
    @{Sheet="Austin"; Cells=@(@{Cell="K2"; Old=196; New=200}, @{Cell="K3"; Old=206; New=208}, @{Cell="K5"; Old=16; New=17}, @{Cell="K6"; Old=226; New=227}, @{Cell="K7"; Old=681; New=689})},
    @{Sheet="South Chicago"; Cells=@(@{Cell="K2"; Old=84; New=85}, @{Cell="K7"; Old=225; New=226})},
    @{Sheet="Garfield Park"; Cells=@(@{Cell="K2"; Old=114; New=115}, @{Cell="K3"; Old=148; New=149}, @{Cell="K6"; Old=118; New=119}, @{Cell="K7"; Old=410; New=413})},
    @{Sheet="Grand Crossing"; Cells=@(@{Cell="K4"; Old=15; New=16}, @{Cell="K6"; Old=106; New=107}, @{Cell="K7"; Old=341; New=343})},
    @{Sheet="New City"; Cells=@(@{Cell="K4"; Old=7; New=8}, @{Cell="K7"; Old=245; New=246})},
    @{Sheet="By Neighborhood"; Cells=@(@{Cell="K2"; Old=81; New=83}, @{Cell="K6"; Old=82; New=83}, @{Cell="K7"; Old=303; New=306}, @{Cell="K8"; Old=681; New=689}, @{Cell="K11"; Old=216; New=220}, @{Cell="K15"; Old=104; New=105}, @{Cell="K16"; Old=33; New=34}, @{Cell="K18"; Old=70; New=73}, @{Cell="K19"; Old=313; New=314}, @{Cell="K20"; Old=239; New=242}, @{Cell="K29"; Old=539; New=541}, @{Cell="K33"; Old=410; New=413}, @{Cell="K36"; Old=121; New=122}, @{Cell="K37"; Old=341; New=343}, @{Cell="K40"; Old=24; New=25}, @{Cell="K41"; Old=88; New=90}, @{Cell="K42"; Old=356; New=358}, @{Cell="K45"; Old=9; New=10}, @{Cell="K48"; Old=120; New=122}, @{Cell="K51"; Old=113; New=114}, @{Cell="J52"; Old=745; New=744}, @{Cell="K52"; Old=285; New=287}, @{Cell="K53"; Old=142; New=144}, @{Cell="K54"; Old=199; New=202}, @{Cell="K57"; Old=28; New=29}, @{Cell="J63"; Old=100; New=101}, @{Cell="K63"; Old=38; New=35}, @{Cell="K65"; Old=245; New=246}, @{Cell="K67"; Old=409; New=412}, @{Cell="K72"; Old=48; New=49}, @{Cell="K73"; Old=93; New=95}, @{Cell="K76"; Old=161; New=162}, @{Cell="K77"; Old=72; New=73}, @{Cell="K79"; Old=265; New=268}, @{Cell="K83"; Old=225; New=226}, @{Cell="K85"; Old=493; New=497}, @{Cell="K86"; Old=67; New=69}, @{Cell="K88"; Old=115; New=116}, @{Cell="K89"; Old=135; New=136}, @{Cell="K91"; Old=105; New=106}, @{Cell="K94"; Old=122; New=123}, @{Cell="K97"; Old=86; New=87}, @{Cell="K101"; Old=10315; New=10388})},
    @{Sheet="North Lawndale"; Cells=@(@{Cell="K2"; Old=126; New=127}, @{Cell="K3"; Old=134; New=135}, @{Cell="K5"; Old=7; New=8}, @{Cell="K7"; Old=409; New=412})},
    @{Sheet="Loop"; Cells=@(@{Cell="K2"; Old=37; New=38}, @{Cell="K6"; Old=89; New=91}, @{Cell="K7"; Old=199; New=202})},
    @{Sheet="Englewood"; Cells=@(@{Cell="K3"; Old=183; New=185}, @{Cell="K7"; Old=539; New=541})},
    @{Sheet="Lake View"; Cells=@(@{Cell="K3"; Old=24; New=25}, @{Cell="K4"; Old=15; New=16}, @{Cell="K7"; Old=120; New=122})},
    @{Sheet="Chatham"; Cells=@(@{Cell="K5"; Old=13; New=14}, @{Cell="K7"; Old=313; New=314})},
    @{Sheet="River North"; Cells=@(@{Cell="K2"; Old=29; New=30}, @{Cell="K7"; Old=161; New=162})},
    @{Sheet="Ashburn"; Cells=@(@{Cell="K2"; Old=29; New=30}, @{Cell="K7"; Old=82; New=83})},
    @{Sheet="Hermosa"; Cells=@(@{Cell="K2"; Old=31; New=32}, @{Cell="K6"; Old=37; New=38}, @{Cell="K7"; Old=88; New=90})},
    @{Sheet="Humboldt Park"; Cells=@(@{Cell="K2"; Old=91; New=92}, @{Cell="K6"; Old=136; New=137}, @{Cell="K7"; Old=356; New=358})},
    @{Sheet="Washington Park"; Cells=@(@{Cell="K2"; Old=27; New=28}, @{Cell="K7"; Old=105; New=106})},
    @{Sheet="Roseland"; Cells=@(@{Cell="K2"; Old=89; New=90}, @{Cell="K3"; Old=94; New=96}, @{Cell="K7"; Old=265; New=268})},
    @{Sheet="Chicago Lawn"; Cells=@(@{Cell="K4"; Old=6; New=7}, @{Cell="K6"; Old=81; New=83}, @{Cell="K7"; Old=239; New=242})},
    @{Sheet="Calumet Heights"; Cells=@(@{Cell="K2"; Old=20; New=21}, @{Cell="K3"; Old=21; New=23}, @{Cell="K7"; Old=70; New=73})},
    @{Sheet="Grand Boulevard"; Cells=@(@{Cell="K4"; Old=10; New=11}, @{Cell="K7"; Old=121; New=122})},
    @{Sheet="Auburn Gresham"; Cells=@(@{Cell="K2"; Old=103; New=104}, @{Cell="K3"; Old=92; New=94}, @{Cell="K7"; Old=303; New=306})},
    @{Sheet="West Loop"; Cells=@(@{Cell="K2"; Old=34; New=35}, @{Cell="K4"; Old=12; New=11}, @{Cell="K6"; Old=51; New=52}, @{Cell="K7"; Old=122; New=123})},
    @{Sheet="Brighton Park"; Cells=@(@{Cell="K2"; Old=35; New=36}, @{Cell="K7"; Old=104; New=105})},
    @{Sheet="Belmont Cragin"; Cells=@(@{Cell="K2"; Old=65; New=66}, @{Cell="K3"; Old=59; New=60}, @{Cell="K4"; Old=9; New=10}, @{Cell="K6"; Old=82; New=83}, @{Cell="K7"; Old=216; New=220})},
    @{Sheet="Portage Park"; Cells=@(@{Cell="K2"; Old=27; New=28}, @{Cell="K3"; Old=22; New=23}, @{Cell="K7"; Old=93; New=95})},
    @{Sheet="Albany Park"; Cells=@(@{Cell="K5"; Old=2; New=3}, @{Cell="K6"; Old=29; New=30}, @{Cell="K7"; Old=81; New=83})},
    @{Sheet="West Town"; Cells=@(@{Cell="K6"; Old=53; New=54}, @{Cell="K7"; Old=86; New=87})},
    @{Sheet="United Center"; Cells=@(@{Cell="K2"; Old=26; New=27}, @{Cell="K7"; Old=115; New=116})},
    @{Sheet="Uptown"; Cells=@(@{Cell="K3"; Old=45; New=46}, @{Cell="K7"; Old=135; New=136})},
    @{Sheet="Streeterville"; Cells=@(@{Cell="K4"; Old=24; New=25}, @{Cell="K6"; Old=15; New=16}, @{Cell="K7"; Old=67; New=69})},
    @{Sheet="Little Italy, UIC"; Cells=@(@{Cell="K3"; Old=28; New=29}, @{Cell="K7"; Old=113; New=114})},
    @{Sheet="Mckinley Park"; Cells=@(@{Cell="K3"; Old=4; New=5}, @{Cell="K7"; Old=28; New=29})},
    @{Sheet="South Shore"; Cells=@(@{Cell="K2"; Old=179; New=181}, @{Cell="K4"; Old=25; New=26}, @{Cell="K6"; Old=110; New=111}, @{Cell="K7"; Old=493; New=497})},
    @{Sheet="Old Town"; Cells=@(@{Cell="K2"; Old=5; New=6}, @{Cell="K7"; Old=48; New=49})},
    @{Sheet="Riverdale"; Cells=@(@{Cell="K3"; Old=25; New=26}, @{Cell="K7"; Old=72; New=73})},
    @{Sheet="Jackson Park"; Cells=@(@{Cell="K3"; Old=2; New=3}, @{Cell="K7"; Old=9; New=10})},
    @{Sheet="Hegewisch"; Cells=@(@{Cell="K2"; Old=8; New=9}, @{Cell="K7"; Old=24; New=25})},
    @{Sheet="Little Village"; Cells=@(@{Cell="K3"; Old=71; New=73}, @{Cell="J4"; Old=30; New=29}, @{Cell="J7"; Old=745; New=744}, @{Cell="K7"; Old=285; New=287})},
    @{Sheet="Bucktown"; Cells=@(@{Cell="K6"; Old=21; New=22}, @{Cell="K7"; Old=33; New=34})}
)

$changed = 0
foreach ($entry in $updates) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    foreach ($c in $entry.Cells) {
        $range = $ws.Range($c.Cell)
        $current = $range.Value2
        if ($current -ne $c.Old) {
            throw "Unexpected existing value in $($entry.Sheet)!$($c.Cell): expected $($c.Old), found $current"
        }
        $range.Value = $c.New
        $changed = $changed + 1
    }
}

Write-Host "Updated $changed cells across $($updates.Count) sheets"